$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) and Volume column (E) stay as text, matching the
# original inline-string cell type (some values look numeric, e.g. "0.691",
# and would otherwise be auto-converted to numbers with float rounding).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "35.130.87"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.904.31"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "253.06"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").Value = "0.691"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "41.49"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("D10").Value = "52.59"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "0.0749"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "0.0981"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "13.06"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Value = "2.181.99"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +3.99%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.939.52"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "4.98"
$ws.Range("E17").Value = "  +3.65%  "
$ws.Range("D18").Value = "35.150.61"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "73.53"
$ws.Range("E19").Value = "  +1.97%  "
$ws.Range("D20").Value = "0.0₃0832"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "243.38"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "12.94"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "5.03"
$ws.Range("E23").Value = "  +5.29%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("E25").Value = "  +5.24%  "
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "166.76"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "8.54"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").Value = "18.52"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "4.128.44"
$ws.Range("E32").Value = "  +12.60%  "
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  -6.88%  "
$ws.Range("D39").Value = "2.00"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "103.46"
$ws.Range("E40").Value = "  +15.75%  "
$ws.Range("D41").Value = "17.21"
$ws.Range("E41").Value = "  +7.27%  "
$ws.Range("D42").Value = "0.0215"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "1.316.27"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "2.43"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("B49").Value = "Gas"
$ws.Range("C49").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D49").Value = "12.21"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "6.58"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").Value = "0.0748"
$ws.Range("E51").Value = "  +6.03%  "
